$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The daily auto-push adds one new log entry for 2026/01/23 in the middle of
# the (date-sorted) table. Insert a new row above the current row 683,
# which shifts the existing rows 683-724 down to 684-725 and grows the used
# range from A1:D724 to A1:D725.
$ws.Rows.Item(683).Insert()

# Write the new entry's day-of-week text and time/ranking numbers directly;
# these are plain values and need no special handling.
$ws.Cells.Item(683, 2).Value = "金"
$ws.Cells.Item(683, 3).Value = 17
$ws.Cells.Item(683, 4).Value = 13

# Column A holds the date as literal text (e.g. "2026/01/23"), matching
# every other row in the sheet. Assigning that string straight to .Value
# would make Excel "smart" auto-convert it into a real date serial number,
# which is not what the source data does. To keep it as plain text without
# leaving a numeric date format behind, stage the text (as Text-formatted)
# in a scratch cell, copy only its value into place, then remove the
# scratch cell entirely so no stray rows/styles remain.
$scratch = $ws.Cells.Item(1000, 1)
$scratch.NumberFormat = "@"
$scratch.Value = "2026/01/23"
$scratch.Copy()
$ws.Cells.Item(683, 1).PasteSpecial(-4163)
$scratch.EntireRow.Delete()
